$d = $word.ActiveDocument

$d.Content.Find.Execute("79×82=", $true, $false, $false, $false, $false, $true, 1, $false, "74×95=", 2)
$d.Content.Find.Execute("98×63=", $true, $false, $false, $false, $false, $true, 1, $false, "34×17=", 2)
$d.Content.Find.Execute("85×58=", $true, $false, $false, $false, $false, $true, 1, $false, "12×50=", 2)
$d.Content.Find.Execute("19×27=", $true, $false, $false, $false, $false, $true, 1, $false, "67×59=", 2)
$d.Content.Find.Execute("35×61=", $true, $false, $false, $false, $false, $true, 1, $false, "72×79=", 2)
$d.Content.Find.Execute("15×59=", $true, $false, $false, $false, $false, $true, 1, $false, "90×66=", 2)
$d.Content.Find.Execute("47×97=", $true, $false, $false, $false, $false, $true, 1, $false, "48×82=", 2)
$d.Content.Find.Execute("51×49=", $true, $false, $false, $false, $false, $true, 1, $false, "23×14=", 2)
$d.Content.Find.Execute("27×30=", $true, $false, $false, $false, $false, $true, 1, $false, "73×58=", 2)
$d.Content.Find.Execute("76×84=", $true, $false, $false, $false, $false, $true, 1, $false, "54×48=", 2)
$d.Content.Find.Execute("60×84=", $true, $false, $false, $false, $false, $true, 1, $false, "20×14=", 2)
$d.Content.Find.Execute("11×82=", $true, $false, $false, $false, $false, $true, 1, $false, "68×63=", 2)
$d.Content.Find.Execute("34×86=", $true, $false, $false, $false, $false, $true, 1, $false, "19×42=", 2)
$d.Content.Find.Execute("17×51=", $true, $false, $false, $false, $false, $true, 1, $false, "15×23=", 2)
$d.Content.Find.Execute("50×94=", $true, $false, $false, $false, $false, $true, 1, $false, "42×38=", 2)
$d.Content.Find.Execute("32×85=", $true, $false, $false, $false, $false, $true, 1, $false, "37×35=", 2)
$d.Content.Find.Execute("63×62=", $true, $false, $false, $false, $false, $true, 1, $false, "90×47=", 2)
$d.Content.Find.Execute("55×96=", $true, $false, $false, $false, $false, $true, 1, $false, "62×47=", 2)
$d.Content.Find.Execute("20×93=", $true, $false, $false, $false, $false, $true, 1, $false, "21×12=", 2)
$d.Content.Find.Execute("31×94=", $true, $false, $false, $false, $false, $true, 1, $false, "99×69=", 2)
$d.Content.Find.Execute("21×47=", $true, $false, $false, $false, $false, $true, 1, $false, "60×16=", 2)
$d.Content.Find.Execute("22×63=", $true, $false, $false, $false, $false, $true, 1, $false, "39×85=", 2)
$d.Content.Find.Execute("60×54=", $true, $false, $false, $false, $false, $true, 1, $false, "72×83=", 2)
$d.Content.Find.Execute("72×82=", $true, $false, $false, $false, $false, $true, 1, $false, "18×67=", 2)
$d.Content.Find.Execute("20×44=", $true, $false, $false, $false, $false, $true, 1, $false, "24×69=", 2)
